$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.827.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.39%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.113.79"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +9.72%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.63"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.51%  "

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.12%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5287"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.05%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4380"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +8.69%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09004"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +7.82%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.05"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +9.10%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.24%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.02"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.05%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.114.90"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +10.25%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.772"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.59%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.811"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.34"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.81%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.19%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001133"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.24%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06667"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.41%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.06"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.90%  "

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.12%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.378"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +7.00%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.903.33"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.61%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.69%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.361.46"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +10.40%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.267"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.68%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.77"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.94%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.565"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +12.78%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.56"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.46"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.41%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.167"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.42%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.36%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.239"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.49%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.024"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.33%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.533"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +21.96%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02613"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.47%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.549"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.22%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06735"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.13%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.530"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +9.76%  "

# Row 40
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.80"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +9.73%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2279"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.83%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6861"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.72%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.254"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.07%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6489"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.91%  "

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.17%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "14.03"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.61%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.247"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.67%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.674"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.36%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.276"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.48%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.58"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.81%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.36"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.55%  "

